$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Agt"
$ws.Cells.Item(2, 3).Value = "Agtr1a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.05804633333333333
$ws.Cells.Item(2, 8).Value = 0.174139
$ws.Cells.Item(2, 9).Value = 0.02760497488446473
$ws.Cells.Item(2, 10).Value = 0.02760497488446473
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.207767333333333
$ws.Cells.Item(2, 14).Value = 6.623302000000001
$ws.Cells.Item(2, 15).Value = 0.07574879669493666
$ws.Cells.Item(2, 16).Value = 0.07574879669493666
$ws.Cells.Item(2, 17).Value = 0.1281527985531111
$ws.Cells.Item(2, 18).Value = 1.153375186978
$ws.Cells.Item(2, 19).Value = 0.002091043630292151
$ws.Cells.Item(2, 20).Value = 0.002091043630292151
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Agt"
$ws.Cells.Item(3, 3).Value = "Agtr1a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.05804633333333333
$ws.Cells.Item(3, 8).Value = 0.174139
$ws.Cells.Item(3, 9).Value = 0.02760497488446473
$ws.Cells.Item(3, 10).Value = 0.02760497488446473
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 18.940215
$ws.Cells.Item(3, 14).Value = 56.820645
$ws.Cells.Item(3, 15).Value = 0.6498413459298955
$ws.Cells.Item(3, 16).Value = 0.6498413459298955
$ws.Cells.Item(3, 17).Value = 1.099410033295
$ws.Cells.Item(3, 18).Value = 9.894690299654998
$ws.Cells.Item(3, 19).Value = 0.01793885403328152
$ws.Cells.Item(3, 20).Value = 0.01793885403328152
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Agt"
$ws.Cells.Item(4, 3).Value = "Agtr1a"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.05804633333333333
$ws.Cells.Item(4, 8).Value = 0.174139
$ws.Cells.Item(4, 9).Value = 0.02760497488446473
$ws.Cells.Item(4, 10).Value = 0.02760497488446473
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.997924
$ws.Cells.Item(4, 14).Value = 23.993772
$ws.Cells.Item(4, 15).Value = 0.2744098573751678
$ws.Cells.Item(4, 16).Value = 0.2744098573751678
$ws.Cells.Item(4, 17).Value = 0.4642501624786667
$ws.Cells.Item(4, 18).Value = 4.178251462307999
$ws.Cells.Item(4, 19).Value = 0.007575077220891057
$ws.Cells.Item(4, 20).Value = 0.007575077220891056
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Agt"
$ws.Cells.Item(5, 3).Value = "Agtr1a"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.761113666666667
$ws.Cells.Item(5, 8).Value = 5.283341
$ws.Cells.Item(5, 9).Value = 0.8375291899635509
$ws.Cells.Item(5, 10).Value = 0.8375291899635509
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.207767333333333
$ws.Cells.Item(5, 14).Value = 6.623302000000001
$ws.Cells.Item(5, 15).Value = 0.07574879669493666
$ws.Cells.Item(5, 16).Value = 0.07574879669493666
$ws.Cells.Item(5, 17).Value = 3.888129223553556
$ws.Cells.Item(5, 18).Value = 34.99316301198201
$ws.Cells.Item(5, 19).Value = 0.063441828336624
$ws.Cells.Item(5, 20).Value = 0.063441828336624
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Agt"
$ws.Cells.Item(6, 3).Value = "Agtr1a"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.761113666666667
$ws.Cells.Item(6, 8).Value = 5.283341
$ws.Cells.Item(6, 9).Value = 0.8375291899635509
$ws.Cells.Item(6, 10).Value = 0.8375291899635509
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 18.940215
$ws.Cells.Item(6, 14).Value = 56.820645
$ws.Cells.Item(6, 15).Value = 0.6498413459298955
$ws.Cells.Item(6, 16).Value = 0.6498413459298955
$ws.Cells.Item(6, 17).Value = 33.35587148610499
$ws.Cells.Item(6, 18).Value = 300.202843374945
$ws.Cells.Item(6, 19).Value = 0.544261096061489
$ws.Cells.Item(6, 20).Value = 0.544261096061489
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Agt"
$ws.Cells.Item(7, 3).Value = "Agtr1a"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.761113666666667
$ws.Cells.Item(7, 8).Value = 5.283341
$ws.Cells.Item(7, 9).Value = 0.8375291899635509
$ws.Cells.Item(7, 10).Value = 0.8375291899635509
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.997924
$ws.Cells.Item(7, 14).Value = 23.993772
$ws.Cells.Item(7, 15).Value = 0.2744098573751678
$ws.Cells.Item(7, 16).Value = 0.2744098573751678
$ws.Cells.Item(7, 17).Value = 14.08525326136133
$ws.Cells.Item(7, 18).Value = 126.767279352252
$ws.Cells.Item(7, 19).Value = 0.2298262655654378
$ws.Cells.Item(7, 20).Value = 0.2298262655654378
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Agt"
$ws.Cells.Item(8, 3).Value = "Agtr1a"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.2331273333333333
$ws.Cells.Item(8, 8).Value = 0.699382
$ws.Cells.Item(8, 9).Value = 0.1108678845327394
$ws.Cells.Item(8, 10).Value = 0.1108678845327394
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.207767333333333
$ws.Cells.Item(8, 14).Value = 6.623302000000001
$ws.Cells.Item(8, 15).Value = 0.07574879669493666
$ws.Cells.Item(8, 16).Value = 0.07574879669493666
$ws.Cells.Item(8, 17).Value = 0.5146909110404444
$ws.Cells.Item(8, 18).Value = 4.632218199364
$ws.Cells.Item(8, 19).Value = 0.008398108845468192
$ws.Cells.Item(8, 20).Value = 0.008398108845468192
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Agt"
$ws.Cells.Item(9, 3).Value = "Agtr1a"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.2331273333333333
$ws.Cells.Item(9, 8).Value = 0.699382
$ws.Cells.Item(9, 9).Value = 0.1108678845327394
$ws.Cells.Item(9, 10).Value = 0.1108678845327394
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 18.940215
$ws.Cells.Item(9, 14).Value = 56.820645
$ws.Cells.Item(9, 15).Value = 0.6498413459298955
$ws.Cells.Item(9, 16).Value = 0.6498413459298955
$ws.Cells.Item(9, 17).Value = 4.41548181571
$ws.Cells.Item(9, 18).Value = 39.73933634138999
$ws.Cells.Item(9, 19).Value = 0.07204653530515565
$ws.Cells.Item(9, 20).Value = 0.07204653530515563
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Agt"
$ws.Cells.Item(10, 3).Value = "Agtr1a"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.2331273333333333
$ws.Cells.Item(10, 8).Value = 0.699382
$ws.Cells.Item(10, 9).Value = 0.1108678845327394
$ws.Cells.Item(10, 10).Value = 0.1108678845327394
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.997924
$ws.Cells.Item(10, 14).Value = 23.993772
$ws.Cells.Item(10, 15).Value = 0.2744098573751678
$ws.Cells.Item(10, 16).Value = 0.2744098573751678
$ws.Cells.Item(10, 17).Value = 1.864534694322667
$ws.Cells.Item(10, 18).Value = 16.780812248904
$ws.Cells.Item(10, 19).Value = 0.0304232403821156
$ws.Cells.Item(10, 20).Value = 0.0304232403821156
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Agt"
$ws.Cells.Item(11, 3).Value = "Agtr1a"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.05046166666666666
$ws.Cells.Item(11, 8).Value = 0.151385
$ws.Cells.Item(11, 9).Value = 0.02399795061924493
$ws.Cells.Item(11, 10).Value = 0.02399795061924493
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.207767333333333
$ws.Cells.Item(11, 14).Value = 6.623302000000001
$ws.Cells.Item(11, 15).Value = 0.07574879669493666
$ws.Cells.Item(11, 16).Value = 0.07574879669493666
$ws.Cells.Item(11, 17).Value = 0.1114076192522222
$ws.Cells.Item(11, 18).Value = 1.00266857327
$ws.Cells.Item(11, 19).Value = 0.001817815882552314
$ws.Cells.Item(11, 20).Value = 0.001817815882552314
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Agt"
$ws.Cells.Item(12, 3).Value = "Agtr1a"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.05046166666666666
$ws.Cells.Item(12, 8).Value = 0.151385
$ws.Cells.Item(12, 9).Value = 0.02399795061924493
$ws.Cells.Item(12, 10).Value = 0.02399795061924493
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 18.940215
$ws.Cells.Item(12, 14).Value = 56.820645
$ws.Cells.Item(12, 15).Value = 0.6498413459298955
$ws.Cells.Item(12, 16).Value = 0.6498413459298955
$ws.Cells.Item(12, 17).Value = 0.9557548159249998
$ws.Cells.Item(12, 18).Value = 8.601793343325
$ws.Cells.Item(12, 19).Value = 0.01559486052996929
$ws.Cells.Item(12, 20).Value = 0.01559486052996929
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Agt"
$ws.Cells.Item(13, 3).Value = "Agtr1a"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.05046166666666666
$ws.Cells.Item(13, 8).Value = 0.151385
$ws.Cells.Item(13, 9).Value = 0.02399795061924493
$ws.Cells.Item(13, 10).Value = 0.02399795061924493
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 7.997924
$ws.Cells.Item(13, 14).Value = 23.993772
$ws.Cells.Item(13, 15).Value = 0.2744098573751678
$ws.Cells.Item(13, 16).Value = 0.2744098573751678
$ws.Cells.Item(13, 17).Value = 0.4035885749133333
$ws.Cells.Item(13, 18).Value = 3.63229717422
$ws.Cells.Item(13, 19).Value = 0.006585274206723322
$ws.Cells.Item(13, 20).Value = 0.006585274206723322
